$wb = $excel.ActiveWorkbook

# --- Rename "MySheet2" to "Sheet2" ---
$sheet2 = $wb.Worksheets.Item("MySheet2")
$sheet2.Name = "Sheet2"

# --- Remove the trailing blank rows 32:37 from Sheet2 ---
$sheet2.Rows("32:37").Delete()

# --- Update the view state: scroll so row 24 is at the top, and move the
#     active selection to J36 (matches the saved sheetView in the target) ---
$sheet2.Activate()
$sheet2.Range("J36").Select()
$excel.ActiveWindow.ScrollRow = 24
$excel.ActiveWindow.ScrollColumn = 1
